$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-7 (years 2004年-2009年), shifting 2010年-2020年 up to rows 2-12
$ws.Range("A2:F7").Delete()

# Row 13 is a brand-new row (2021年) - copy the existing label style from A12
# so the new A13 cell picks up the same formatting (bold, border, centered).
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 22062
$ws.Range("C13").Value = 64.7962
$ws.Range("D13").Value = 83.5659
$ws.Range("E13").Value = 347.9788
$ws.Range("F13").Value = 42.42
